$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Q0)
$ws.Range("B3").Value = 0.0944151529069
$ws.Range("C3").Value = 0.4355155406315797
$ws.Range("D3").Value = 0.2568765859223691
$ws.Range("E3").Value = 0.5068299378710467
$ws.Range("F3").Value = 0.5154357014326421
$ws.Range("G3").Value = 15

# Row 4 (Q1)
$ws.Range("B4").Value = 0.3384128879484091
$ws.Range("C4").Value = 0.586491187861096
$ws.Range("D4").Value = 0.5772465913568255
$ws.Range("E4").Value = 0.7597674587377545
$ws.Range("F4").Value = 0.705916007249359
$ws.Range("G4").Value = 14
